$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Databases")

# Add new "Role" header in column X, copying the formatting of W1 (bold/centered header)
$ws.Range("X1").Value = "Role"
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the view: scroll so column N is the top-left visible column, and select X1
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("X1").Select()
